# Parametrizacion de cedulas con cierre de sesion - 10/08/2025
# Appends a new batch of cedula rows (date marker + IDs with validation-status fills)
# to column A of Hoja1, rows 1764:1905.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Write the new cedula values (and the leading date marker) in one batch.
$arr = New-Object 'object[,]' 142,1
$arr[0,0] = 45938
$arr[1,0] = 21205258
$arr[2,0] = 31006459
$arr[3,0] = 11373533
$arr[4,0] = 30030456
$arr[5,0] = 70052315
$arr[6,0] = 40449322
$arr[7,0] = 40329872
$arr[8,0] = 40219846
$arr[9,0] = 1122338100
$arr[10,0] = 17303753
$arr[11,0] = 17345135
$arr[12,0] = 33875522
$arr[13,0] = 36293795
$arr[14,0] = 20871642
$arr[15,0] = 40186554
$arr[16,0] = 52287828
$arr[17,0] = 40333471
$arr[18,0] = 50976844
$arr[19,0] = 40329979
$arr[20,0] = 35263486
$arr[21,0] = 40415999
$arr[22,0] = 11405158
$arr[23,0] = 40384750
$arr[24,0] = 26508842
$arr[25,0] = 1121920935
$arr[26,0] = 30081591
$arr[27,0] = 40371629
$arr[28,0] = 1121840691
$arr[29,0] = 41213156
$arr[30,0] = 30080583
$arr[31,0] = 52812191
$arr[32,0] = 17329021
$arr[33,0] = 39728361
$arr[34,0] = 21011250
$arr[35,0] = 1120868086
$arr[36,0] = 21243208
$arr[37,0] = 1799801717
$arr[38,0] = 15960326
$arr[39,0] = 1120868173
$arr[40,0] = 900084777
$arr[41,0] = 1033690960
$arr[42,0] = 1075315269
$arr[43,0] = 33676007
$arr[44,0] = 1121818890
$arr[45,0] = 555555555
$arr[46,0] = 555555556
$arr[47,0] = 900737989
$arr[48,0] = 9007379890
$arr[49,0] = 1121913651
$arr[50,0] = 9000847773
$arr[51,0] = 9000847774
$arr[52,0] = 52861841
$arr[53,0] = 555555551
$arr[54,0] = 5849675
$arr[55,0] = 555555552
$arr[56,0] = 555555553
$arr[57,0] = 555555554
$arr[58,0] = 40417192
$arr[59,0] = 1122117543
$arr[60,0] = 1799831717
$arr[61,0] = 1799821717
$arr[62,0] = 1799841717
$arr[63,0] = 40416523
$arr[64,0] = 9000847771
$arr[65,0] = 9000847772
$arr[66,0] = 40188434
$arr[67,0] = 1121828917
$arr[68,0] = 1006771687
$arr[69,0] = 1118555784
$arr[70,0] = 79891869
$arr[71,0] = 40438899
$arr[72,0] = 555555557
$arr[73,0] = 68297670
$arr[74,0] = 14952515
$arr[75,0] = 1033723546
$arr[76,0] = 1120374995
$arr[77,0] = 86044711
$arr[78,0] = 40342629
$arr[79,0] = 1006720164
$arr[80,0] = 1121869388
$arr[81,0] = 890900608214
$arr[82,0] = 890900608211
$arr[83,0] = 555555558
$arr[84,0] = 89090060825
$arr[85,0] = 40356371
$arr[86,0] = 40218334
$arr[87,0] = 1108928812
$arr[88,0] = 9007505342
$arr[89,0] = 1003483201
$arr[90,0] = 1121911590
$arr[91,0] = 30080723
$arr[92,0] = 1006779023
$arr[93,0] = 8909006082
$arr[94,0] = 890900608213
$arr[95,0] = 89090060821
$arr[96,0] = 89090060829
$arr[97,0] = 890900608212
$arr[98,0] = 89090060824
$arr[99,0] = 89090060827
$arr[100,0] = 89090060828
$arr[101,0] = 890900608210
$arr[102,0] = 40305749
$arr[103,0] = 1006718970
$arr[104,0] = 1122919735
$arr[105,0] = 1006859801
$arr[106,0] = 1004634852
$arr[107,0] = 1116809197
$arr[108,0] = 444444444
$arr[109,0] = 40330136
$arr[110,0] = 1038100731
$arr[111,0] = 1116861125
$arr[112,0] = 1120866044
$arr[113,0] = 69022144
$arr[114,0] = 1121820427
$arr[115,0] = 1053764354
$arr[116,0] = 1011086251
$arr[117,0] = 1121849836
$arr[118,0] = 1120864873
$arr[119,0] = 1133839243
$arr[120,0] = 1121968001
$arr[121,0] = 77777778
$arr[122,0] = 77777779
$arr[123,0] = 1006719384
$arr[124,0] = 1006858033
$arr[125,0] = 1122922330
$arr[126,0] = 1082773285
$arr[127,0] = 55190862
$arr[128,0] = 1070325222
$arr[129,0] = 30520140
$arr[130,0] = 1124818008
$arr[131,0] = 444444441
$arr[132,0] = 444444442
$arr[133,0] = 444444443
$arr[134,0] = 1006719691
$arr[135,0] = 1121946622
$arr[136,0] = 1006720627
$arr[137,0] = 1123565052
$arr[138,0] = 40413517
$arr[139,0] = 63251940
$arr[140,0] = 1006009495
$arr[141,0] = 1002596538
$ws.Range("A1764:A1905").Value = $arr

# 2) Re-apply the existing per-row validation-status formatting by copying
#    format from a representative cell already using that style, so the
#    workbook keeps reusing the same style indices instead of minting new ones.

# style s="1" (date marker (numFmt date))
$rngList_1 = @("A1764")
foreach ($rng in $rngList_1) {
    $ws.Range("A1").Copy() | Out-Null
    $ws.Range($rng).PasteSpecial(-4122) | Out-Null
}

# style s="3" (green fill - validated cedula)
$rngList_3 = @("A1766", "A1770:A1773", "A1776:A1779", "A1789:A1793", "A1798:A1805", "A1807:A1809", "A1811:A1816", "A1818", "A1822:A1826", "A1828:A1830", "A1832:A1836", "A1839:A1844", "A1847", "A1849:A1856", "A1866:A1888", "A1890:A1905")
foreach ($rng in $rngList_3) {
    $ws.Range("A13").Copy() | Out-Null
    $ws.Range($rng).PasteSpecial(-4122) | Out-Null
}

# style s="5" (red fill - rejected/duplicate cedula)
$rngList_5 = @("A1765", "A1767:A1769", "A1774:A1775", "A1786:A1788", "A1794:A1797", "A1806", "A1810", "A1817", "A1819:A1821", "A1827", "A1831", "A1838", "A1845:A1846", "A1848", "A1857:A1865", "A1889")
foreach ($rng in $rngList_5) {
    $ws.Range("A39").Copy() | Out-Null
    $ws.Range($rng).PasteSpecial(-4122) | Out-Null
}

# style s="7" (green fill + explicit font color)
$rngList_7 = @("A1780:A1784")
foreach ($rng in $rngList_7) {
    $ws.Range("A1173").Copy() | Out-Null
    $ws.Range($rng).PasteSpecial(-4122) | Out-Null
}

# style s="8" (red fill + explicit font color)
$rngList_8 = @("A1785")
foreach ($rng in $rngList_8) {
    $ws.Range("A1178").Copy() | Out-Null
    $ws.Range($rng).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false

# 3) Match the author's final scroll position / selection.
$ws.Range("A1726").Select() | Out-Null
$ws.Range("B1743").Select() | Out-Null
